$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Balance Sheet")
$ws2 = $wb.Worksheets.Item("Income Statement")

function Set-TextValue($ws, $addr, $text, $scratch) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# Scratch cell used to force numeric-looking strings to be stored as literal text
$scratch = $ws1.Range("Z100")
$scratch.NumberFormat = "@"

# ---- sheet1 ----
$ws1.Range("A1").Value = ' **Balance Sheet Indicators**'
$ws1.Range("B1").Value = '**Current Year (2024)**'
$ws1.Range("C1").Value = '**Previous Year (2023)** '
$ws1.Range("A2").Value = ' Cash and Cash Equivalents'
Set-TextValue $ws1 "B2" '1000' $scratch
Set-TextValue $ws1 "C2" '1000                     ' $scratch
$ws1.Range("A3").Value = ' Accounts Receivable'
Set-TextValue $ws1 "B3" '11,987,605.97' $scratch
Set-TextValue $ws1 "C3" '10,711,454.12            ' $scratch
$ws1.Range("A4").Value = ' Property, Plant and Equipment (Net)'
Set-TextValue $ws1 "B4" '3,494,523.92' $scratch
Set-TextValue $ws1 "C4" '3,494,523.92             ' $scratch
$ws1.Range("A5").Value = ' Total Assets'
Set-TextValue $ws1 "B5" '14,355,193.96' $scratch
Set-TextValue $ws1 "C5" '13,424,369.47            ' $scratch
$ws1.Range("A6").Value = ' Accounts Payable'
Set-TextValue $ws1 "B6" '-12,443,892.15' $scratch
Set-TextValue $ws1 "C6" '-10,979,515.78           ' $scratch
$ws1.Range("A7").Value = ' Accumulated Profit/(Loss)'
Set-TextValue $ws1 "B7" '-2,444,853.69' $scratch
Set-TextValue $ws1 "C7" '-2,741,596.38            ' $scratch
$ws1.Range("A8").Value = ' Total Liabilities'
Set-TextValue $ws1 "B8" '-14,888,745.84' $scratch
Set-TextValue $ws1 "C8" '-13,721,112.16           ' $scratch

# ---- sheet2 ----
$ws2.Range("A1").Value = ' **Income Statement Indicators**'
$ws2.Range("B1").Value = '**Current Year (2024)**'
$ws2.Range("C1").Value = '**Previous Year (2023)** '
$ws2.Range("A2").Value = ' Revenue'
Set-TextValue $ws2 "B2" '-1,276,151.85' $scratch
Set-TextValue $ws2 "C2" '-1,727,145.61            ' $scratch
$ws2.Range("A3").Value = ' Cost of Goods Sold'
Set-TextValue $ws2 "B3" '367,148.33' $scratch
Set-TextValue $ws2 "C3" '428,073.69               ' $scratch
$ws2.Range("A4").Value = ' Gross Profit'
$ws2.Range("B4").Value = '*N/A*'
$ws2.Range("C4").Value = '*N/A*                    '
$ws2.Range("A5").Value = ' General and Administrative Expenses'
$ws2.Range("B5").Value = '*Consolidated within Expense Total*'
$ws2.Range("C5").Value = '*Consolidated within Expense Total* '
$ws2.Range("A6").Value = ' Net Profit'
$ws2.Range("B6").Value = '*N/A*'
$ws2.Range("C6").Value = '*N/A*                    '

# ---- Add Cash Flow Statement sheet ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Cash Flow Statement"

$ws3.Range("A1").Value = ' **Cash Flow Indicators**'
$ws3.Range("B1").Value = '**Current Year (2024)**'
$ws3.Range("C1").Value = '**Previous Year (2023)** '
$ws3.Range("A2").Value = ' Depreciation and Amortization'
Set-TextValue $ws3 "B2" '350,277.36' $scratch
Set-TextValue $ws3 "C2" '349,452.36               ' $scratch
$ws3.Range("A3").Value = ' Total Cash Flow'
Set-TextValue $ws3 "B3" '533,551.88' $scratch
Set-TextValue $ws3 "C3" '296,742.69               ' $scratch

# Copy header formatting (bold, border, centered) from an existing header row
$ws1.Range("A1:C1").Copy()
$ws3.Range("A1:C1").PasteSpecial(-4122)

# Clean up scratch cell entirely so it leaves no trace in the sheet
$scratch.EntireColumn.Delete()
